# FAST_holdings.xlsx update:
#  - bump the "as of" date in the confidentiality banner (A13) from 2021-07-09 to 2021-07-13
#  - refresh the Weight / Percent Change figures in D2:E10 for the new date
#
# The worksheet ships sheet-protected, so it must be unprotected before the
# cells can be written, then protection is restored afterward.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Unprotect()

$ws.Range("A13").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-07-13 for illustrative purposes only and are subject to change."

$ws.Range("D2").Value = 0.1347410100502228
$ws.Range("E2").Value = -0.02395480225988711

$ws.Range("D3").Value = 0.1096951448586402
$ws.Range("E3").Value = -0.01030219780219788

$ws.Range("D4").Value = 0.1128694732263445
$ws.Range("E4").Value = -0.001842570754716943

$ws.Range("D5").Value = 0.1188656012081007
$ws.Range("E5").Value = -0.009343647085854334

$ws.Range("D6").Value = 0.1212665951282338
$ws.Range("E6").Value = -0.002133048926809988

$ws.Range("D7").Value = 0.1434908176642477
$ws.Range("E7").Value = -0.009706889988580159

$ws.Range("D8").Value = 0.1312155655552166
$ws.Range("E8").Value = -0.003592152528322567

$ws.Range("D9").Value = 0.1278557923089938
$ws.Range("E9").Value = -0.009542385654234486

$ws.Range("D10").Value = 1
$ws.Range("E10").Value = -0.009019316317067427

$ws.Protect()
